$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3225525858620018
$ws.Range("D2").Value = 0.7500794634550545

$ws.Range("C3").Value = -0.4021485301734843
$ws.Range("D3").Value = 0.6914551560940825

$ws.Range("C4").Value = -2.12083054957058
$ws.Range("D4").Value = 0.04544542231610049

$ws.Range("C5").Value = -2.209759277230824
$ws.Range("D5").Value = 0.03783166011216088

$ws.Range("C6").Value = -0.2154264278911612
$ws.Range("D6").Value = 0.8314195293148137

$ws.Range("C7").Value = -1.654450946312143
$ws.Range("D7").Value = 0.1122344321331386

$ws.Range("C8").Value = -2.234511832715092
$ws.Range("D8").Value = 0.03592908887941393

$ws.Range("C9").Value = -1.840782545278833
$ws.Range("D9").Value = 0.07918127762550098

$ws.Range("C10").Value = -2.384512789139425
$ws.Range("D10").Value = 0.02615155055424601

$ws.Range("C11").Value = -0.5937314598320715
$ws.Range("D11").Value = 0.5587476843922476
